# This workbook is a weekly/daily price log (Acelga, Feria Lagunitas de
# Puerto Montt). The edit inserts two new observation rows into the table
# (one after the existing row 43, one further down in the list), which
# pushes every subsequent row down by one. Column order in every row is
# fixed: A Mercado ID, B Mercado, C Region, D Fecha, E Codreg,
# F Categoria ID, G Categoria, H Variedad, I Calidad, J Volumen,
# K Precio minimo, L Precio maximo, M Precio promedio ponderado,
# N Unidad de comercializacion, O Origen, P Precio $/Kg, Q Kg o Unidades,
# R Clasificacion. Columns A,B,C,E,F,G,H,I,R are constant for every data
# row in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$constA = 4
$constB = "Feria Lagunitas de Puerto Montt"
$constC = "Los Lagos"
$constE = 10
$constF = 100112009
$constG = "Acelga"
$constH = "Sin especificar"
$constI = "Primera"
$constR = "Hortaliza"

function Set-DataRow {
    param(
        [int]$RowNum,
        [double]$D,
        [double]$J,
        [double]$K,
        [double]$L,
        [double]$M,
        [string]$N,
        [string]$O,
        [double]$P,
        [double]$Q
    )
    $ws.Cells.Item($RowNum, 1).Value2 = $constA
    $ws.Cells.Item($RowNum, 2).Value2 = $constB
    $ws.Cells.Item($RowNum, 3).Value2 = $constC
    $ws.Cells.Item($RowNum, 4).Value2 = $D
    $ws.Cells.Item($RowNum, 5).Value2 = $constE
    $ws.Cells.Item($RowNum, 6).Value2 = $constF
    $ws.Cells.Item($RowNum, 7).Value2 = $constG
    $ws.Cells.Item($RowNum, 8).Value2 = $constH
    $ws.Cells.Item($RowNum, 9).Value2 = $constI
    $ws.Cells.Item($RowNum, 10).Value2 = $J
    $ws.Cells.Item($RowNum, 11).Value2 = $K
    $ws.Cells.Item($RowNum, 12).Value2 = $L
    $ws.Cells.Item($RowNum, 13).Value2 = $M
    $ws.Cells.Item($RowNum, 14).Value2 = $N
    $ws.Cells.Item($RowNum, 15).Value2 = $O
    $ws.Cells.Item($RowNum, 16).Value2 = $P
    $ws.Cells.Item($RowNum, 17).Value2 = $Q
    $ws.Cells.Item($RowNum, 18).Value2 = $constR
}

# Insert a new row at 44 (old rows 44..147 shift down to 45..148),
# then populate the new row with the new observation.
$ws.Rows.Item(44).Insert()
Set-DataRow 44 44614 80 10000 10000 10000 "`$/docena de atados (12 kilos)" "Región de La Araucanía" 833 12

# Insert a second new row at 141 (current rows 141..148 shift down to
# 142..149), then populate the new row with the new observation.
$ws.Rows.Item(141).Insert()
Set-DataRow 141 44615 20 10000 10000 10000 "`$/docena de atados (12 kilos)" "Región de La Araucanía" 833 12
